$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Locate the 1-based Paragraphs index of the paragraph that contains the
# (unique) search text $needle.
function Get-ParaIndexByText([string]$needle) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw "text not found: $needle" }
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            return $i
        }
    }
    throw "paragraph index not found for: $needle"
}

# Insert a brand-new, format-free list paragraph (same pStyle/numPr as the
# clean template paragraph "admin felulet...") immediately before the
# paragraph currently at 1-based index $idx, fill it with $newText, and
# (optionally) color its text/paragraph-mark with $colorVal (a wdColor
# integer). Returns nothing; the paragraph previously at $idx is now at
# $idx + 1 (the caller is responsible for deleting it afterwards if needed).
function New-CleanListParagraphBefore([int]$idx, [string]$newText, $colorVal) {
    $templateIdx = Get-ParaIndexByText("admin")
    $template = $d.Paragraphs.Item($templateIdx)
    $dup = $template.Range.Duplicate

    $target = $d.Paragraphs.Item($idx)
    $insertPos = $target.Range.Start

    $ins = $d.Range($insertPos, $insertPos)
    $ins.FormattedText = $dup.FormattedText

    # The freshly-inserted (still colorless) paragraph is now at $idx.
    $newPara = $d.Paragraphs.Item($idx)
    $newFull = $newPara.Range
    $newTextOnly = $d.Range($newFull.Start, $newFull.End - 1)
    $newTextOnly.Text = $newText

    if ($colorVal -ne $null) {
        $coloredPara = $d.Paragraphs.Item($idx)
        $coloredPara.Range.Font.Color = $colorVal
    }
}

# Replace the paragraph at 1-based index $idx (identified earlier by
# $needle, re-resolved fresh to stay correct) with a colorless copy of the
# same list paragraph containing $newText - i.e. strip any direct color
# formatting that paragraph/run carried.
function Set-ParagraphTextNoColor([string]$needle, [string]$newText) {
    $idx = Get-ParaIndexByText($needle)
    New-CleanListParagraphBefore $idx $newText $null
    # Original (now shifted down by one) paragraph follows right after;
    # delete it, mark included, which merges everything back cleanly.
    $old = $d.Paragraphs.Item($idx + 1)
    $old.Range.Delete()
}

# ---------------------------------------------------------------------------
# 1) Strip the (now unwanted) direct color formatting from three bullets.
# ---------------------------------------------------------------------------

Set-ParagraphTextNoColor "Bejelentkez" `
    "Bejelentkezésnél először kelljen megnyomni a gombot és utána kelljen beírni a jelszót és a felhasználót"

Set-ParagraphTextNoColor "Téma hozzáadásánál" `
    "Téma hozzáadásánál is kelljen előbb gombot megnyomni és csak utána ugorjon fel chat ablak ahová beírhatod az új témát"

Set-ParagraphTextNoColor "Főmenü gomb" `
    "Főmenü gomb"

# ---------------------------------------------------------------------------
# 2) Split the "Sötét mód fixálása" bullet into two new bullets:
#      - a new red ("FF0000") bullet about deleting cases
#      - a plain bullet "Sötét mód hozzáadása" (replacing the old wording)
# ---------------------------------------------------------------------------

$oldIdx = Get-ParaIndexByText("fixálása")

New-CleanListParagraphBefore $oldIdx "Ügyintézés alatt és lezárt ügyek törlése" 255
New-CleanListParagraphBefore ($oldIdx + 1) "Sötét mód hozzáadása" $null

# Delete the original "Sötét mód fixálása" paragraph, now shifted two slots
# further down by the two insertions above.
$old = $d.Paragraphs.Item($oldIdx + 2)
$old.Range.Delete()

Write-Output "done"
